$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Set column B header and definitions for each location type row
$ws.Range("B1").Value = 'definition'
$ws.Range("B2").Value = 'A sampling location within a facility that is not better described by any other location type.'
$ws.Range("B3").Value = 'A natural or anthropogenic body of fresh water, not used for wastewater storage or disposal.'
$ws.Range("B4").Value = 'An anthropogenic body of water designed to disperse into the underlying soil, typically for pollution abatement.'
$ws.Range("B5").Value = 'A sewer system that carries a mixture of stormwater and sanitary sewage.'
$ws.Range("B6").Value = 'The ocean or open sea.'
$ws.Range("B7").Value = 'A land-based sampling location for soil, soil vapor, plants, etc.'
$ws.Range("B8").Value = 'A location where groundwater naturally emerges at the surface.'
$ws.Range("B9").Value = 'Any type of groundwater well.'
$ws.Range("B10").Value = 'A sewer system that carries stormwater only, not combined or sanitary sewage.'
$ws.Range("B11").Value = 'A brackish water zone where a river meets the ocean.'
$ws.Range("B12").Value = 'Water seeping from the ground, sometimes forming small pools.'
$ws.Range("B13").Value = 'A sewer system that carries only sanitary (wastewater) flow, not combined with stormwater.'
$ws.Range("B14").Value = 'A site containing samples of refuse or related sediments.'
$ws.Range("B15").Value = 'A naturally flowing body of fresh water.'
$ws.Range("B16").Value = 'A point where air emissions are sampled near or at their release into the atmosphere.'
$ws.Range("B17").Value = 'An anthropogenic trench or pipe used to drain stormwater.'
$ws.Range("B18").Value = 'Similar to a terrestrial location, but specifically used for agricultural purposes.'
$ws.Range("B19").Value = 'A location at or very near where wastewater is discharged into the environment.'
$ws.Range("B20").Value = 'An area with buried perforated pipes used to drain saturated soil.'
$ws.Range("B21").Value = 'A site established to monitor air quality and related parameters.'
$ws.Range("B22").Value = 'Water or solids applied to land for irrigation or disposal.'
$ws.Range("B23").Value = 'An anthropogenic pond designed to slow and temporarily hold runoff, typically for pollution control.'
$ws.Range("B24").Value = 'A sampling point near a storage facility, often for industrial or waste monitoring.'
$ws.Range("B25").Value = 'A well used to inject wastewater into an aquifer.'
$ws.Range("B26").Value = 'The burning of waste or vegetation (slash) in the open air.'
$ws.Range("B27").Value = 'A tank used to separate and collect solids from liquid household wastewater.'

# Adjust column widths to match the final layout
$ws.Columns.Item(1).ColumnWidth = 30
$ws.Columns.Item(2).ColumnWidth = 128

# Restore selection state
$ws.Range("B9").Select()
